# Add "NA" values under the duplicate_image_filename column (column E)
# for the rows that currently contain data (rows 2 through 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
